# Edit the "Dashboard" slide layout (slideLayout12.xml) inside the slide master:
#  1. Remove the "FINANCES" TextBox (cNvPr id=10, name="TextBox 9").
#  2. Move the horizontal divider line (Straight Connector 159, id=160) down.
#  3. Resize/reposition the four dashed placeholder rectangles
#     (Rectangle 11/12/13/14, ids 12/13/14/15) to be taller, giving more
#     room for the plots/text underneath them.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layout = $master.CustomLayouts.Item(12)

$shapes = $layout.Shapes

# 1. Delete the "FINANCES" textbox (first shape in the tree).
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $sh = $shapes.Item($i)
    if ($sh.Name -eq "TextBox 9") {
        $sh.Delete()
    }
}

# 2. Move the straight connector line down from y=4688544 EMU to y=5029200 EMU.
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -eq "Straight Connector 159") {
        $sh.Top = 396.0
    }
}

# 3. Resize the four dashed rectangles: new top=1389078 EMU, new height=685800 EMU
#    (was top=1463040 EMU, height=274320 EMU).
$rectNames = @("Rectangle 11", "Rectangle 12", "Rectangle 13", "Rectangle 14")
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($rectNames -contains $sh.Name) {
        $sh.Top = 109.37622047244095
        $sh.Height = 54.0
    }
}
